$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("LoginApp")
$wsPayment = $wb.Worksheets.Item("PaymentPage")

# New test-case row appended to the LoginApp sheet.
$wsLogin.Range("B18").Value = 'Check wheather Torenzocafe Title is displayed while clicking on "Cancel" from Role page.'
$wsLogin.Range("C18").Value = "PASS"

# PaymentPage was the selected/active tab before; keep its selection as-is,
# but do this BEFORE activating LoginApp so the later Activate() sticks.
$wsPayment.Activate()
$wsPayment.Range("A14").Select()

# LoginApp becomes the active tab, with the new row's first empty cell selected.
$wsLogin.Activate()
$wsLogin.Range("B19").Select()
